# Update "想去人数" (column F, interest count) figures across the
# workbook's four sheets to match the refreshed gh-pages data export.
#
# Sheet order (per xl/workbook.xml): 1=展览, 2=演出, 3=本地生活, 4=全部类型
# (全部类型 is a concatenation of the other three sheets' rows, hence the
# repeated values across sheets).

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(4, 6).Value = 5881
$ws.Cells.Item(5, 6).Value = 72
$ws.Cells.Item(6, 6).Value = 62
$ws.Cells.Item(9, 6).Value = 1568
$ws.Cells.Item(11, 6).Value = 30
$ws.Cells.Item(13, 6).Value = 1580
$ws.Cells.Item(14, 6).Value = 1580
$ws.Cells.Item(15, 6).Value = 1538
$ws.Cells.Item(17, 6).Value = 145
$ws.Cells.Item(18, 6).Value = 613
$ws.Cells.Item(19, 6).Value = 4401
$ws.Cells.Item(22, 6).Value = 3335
$ws.Cells.Item(23, 6).Value = 813
$ws.Cells.Item(24, 6).Value = 4
$ws.Cells.Item(26, 6).Value = 2301
$ws.Cells.Item(30, 6).Value = 451
$ws.Cells.Item(32, 6).Value = 786
$ws.Cells.Item(34, 6).Value = 1198
$ws.Cells.Item(35, 6).Value = 1186

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(12, 6).Value = 105
$ws.Cells.Item(15, 6).Value = 37
$ws.Cells.Item(18, 6).Value = 125
$ws.Cells.Item(19, 6).Value = 298
$ws.Cells.Item(20, 6).Value = 227
$ws.Cells.Item(21, 6).Value = 490

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 6).Value = 636
$ws.Cells.Item(4, 6).Value = 176
$ws.Cells.Item(5, 6).Value = 261

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(6, 6).Value = 636
$ws.Cells.Item(7, 6).Value = 176
$ws.Cells.Item(8, 6).Value = 5881
$ws.Cells.Item(10, 6).Value = 72
$ws.Cells.Item(11, 6).Value = 62
$ws.Cells.Item(20, 6).Value = 1568
$ws.Cells.Item(23, 6).Value = 30
$ws.Cells.Item(24, 6).Value = 1580
$ws.Cells.Item(25, 6).Value = 105
$ws.Cells.Item(26, 6).Value = 1538
$ws.Cells.Item(28, 6).Value = 145
$ws.Cells.Item(29, 6).Value = 613
$ws.Cells.Item(30, 6).Value = 4401
$ws.Cells.Item(32, 6).Value = 3335
$ws.Cells.Item(33, 6).Value = 813
$ws.Cells.Item(35, 6).Value = 2301
$ws.Cells.Item(39, 6).Value = 451
$ws.Cells.Item(42, 6).Value = 125
$ws.Cells.Item(43, 6).Value = 298
$ws.Cells.Item(44, 6).Value = 227
$ws.Cells.Item(45, 6).Value = 490
$ws.Cells.Item(46, 6).Value = 786
